$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.257483243942261
$ws.Range("B1").Value = 2.56754732131958
$ws.Range("C1").Value = 4.931248664855957
$ws.Range("D1").Value = 2.003765106201172
$ws.Range("E1").Value = 1.154708027839661
